$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.777.69"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.861.31"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.58%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.31"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.614"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.90%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.70"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +11.77%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +6.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0699"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.08%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.33%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.48%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.858.16"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.73%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +7.38%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +7.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.712.16"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.49%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.39"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0806"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "246.42"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.28"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +8.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.80"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +15.72%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.33"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.03"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.23%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.17%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.46"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +18.20%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.306.61"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +36.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.95"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0545"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.59%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.91"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.23%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "Aave"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "93.65"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +14.36%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.689"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +6.14%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.52"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +6.94%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.345.97"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.70%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.03"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +7.93%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0196"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.58%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "15.20"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +9.50%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.96%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +2.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.83"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0521"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +2.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.11"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.026.37"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "104.79"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.17%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.26%  "
